$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 730.6667
$ws.Range("I32").Value = 620
$ws.Range("J32").Value = 786
$ws.Range("K32").Value = 620
$ws.Range("L32").Value = 786
$ws.Range("M32").Value = -294
$ws.Range("N32").Value = -1438
$ws.Range("H51").Value = 5828.7144
$ws.Range("I51").Value = 7749.75
$ws.Range("J51").Value = 3267.3333
$ws.Range("K51").Value = 7749.75
$ws.Range("L51").Value = 3267.3333
$ws.Range("M51").Value = -7265.75
$ws.Range("N51").Value = -4235.3333
$ws.Range("H112").Value = 6383937.5
$ws.Range("J112").Value = 3177297
$ws.Range("L112").Value = 9531891
$ws.Range("N112").Value = -9534107
$ws.Range("H127").Value = 1432.6666
$ws.Range("J127").Value = 1732.3334
$ws.Range("L127").Value = 5197.0002
$ws.Range("N127").Value = -15117.0002
$ws.Range("H129").Value = 400974.38
$ws.Range("J129").Value = 589530.1
$ws.Range("L129").Value = 1768590.3
$ws.Range("N129").Value = -1778590.3
$ws.Range("H138").Value = 4446.478
$ws.Range("I138").Value = 3784.7144
$ws.Range("J138").Value = 4565.2563
$ws.Range("K138").Value = 11354.1432
$ws.Range("L138").Value = 13695.7689
$ws.Range("M138").Value = -6214.143199999999
$ws.Range("N138").Value = -23975.7689

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8785654
$ws.Range("I61").Value = 10335378
$ws.Range("K61").Value = 10335378
$ws.Range("M61").Value = -10335166
$ws.Range("H74").Value = 41670540
$ws.Range("I74").Value = 90913540
$ws.Range("K74").Value = 90913540
$ws.Range("M74").Value = -90912666
$ws.Range("H77").Value = 41670540
$ws.Range("I77").Value = 90913540
$ws.Range("K77").Value = 454567700
$ws.Range("M77").Value = -454563332
$ws.Range("H122").Value = 1746.2941
$ws.Range("I122").Value = 1813.8
$ws.Range("J122").Value = 1240
$ws.Range("K122").Value = 5441.4
$ws.Range("L122").Value = 3720
$ws.Range("M122").Value = -2991.4
$ws.Range("N122").Value = -8620
$ws.Range("H136").Value = 8785654
$ws.Range("I136").Value = 10335378
$ws.Range("K136").Value = 31006134
$ws.Range("M136").Value = -31003584
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6993.875
$ws.Range("J31").Value = 9708.25
$ws.Range("L31").Value = 9708.25
$ws.Range("N31").Value = -10298.25
$ws.Range("H34").Value = 6993.875
$ws.Range("J34").Value = 9708.25
$ws.Range("L34").Value = 9708.25
$ws.Range("N34").Value = -10112.25
$ws.Range("H58").Value = 27092
$ws.Range("I58").Value = 2492.2856
$ws.Range("J58").Value = 40338
$ws.Range("K58").Value = 2492.2856
$ws.Range("L58").Value = 40338
$ws.Range("M58").Value = -2289.2856
$ws.Range("N58").Value = -40744
$ws.Range("H86").Value = 8958.471
$ws.Range("I86").Value = 2165.125
$ws.Range("J86").Value = 14997
$ws.Range("K86").Value = 2165.125
$ws.Range("L86").Value = 14997
$ws.Range("M86").Value = -1042.125
$ws.Range("N86").Value = -17243
$ws.Range("H89").Value = 8958.471
$ws.Range("I89").Value = 2165.125
$ws.Range("J89").Value = 14997
$ws.Range("K89").Value = 10825.625
$ws.Range("L89").Value = 74985
$ws.Range("M89").Value = -5209.625
$ws.Range("N89").Value = -86217
$ws.Range("H136").Value = 27092
$ws.Range("I136").Value = 2492.2856
$ws.Range("J136").Value = 40338
$ws.Range("K136").Value = 7476.8568
$ws.Range("L136").Value = 121014
$ws.Range("M136").Value = -4926.8568
$ws.Range("N136").Value = -126114
$ws.Range("H141").Value = 37497
$ws.Range("J141").Value = 39996.727
$ws.Range("L141").Value = 39996.727
$ws.Range("N141").Value = -50356.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H131").Value = 768.8200000000001
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 768.8200000000001
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2306.46
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12386.46

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4596
$ws.Range("I97").Value = 3388.5
$ws.Range("J97").Value = 7011
$ws.Range("K97").Value = 3388.5
$ws.Range("L97").Value = 7011
$ws.Range("M97").Value = -2892.5
$ws.Range("N97").Value = -8003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2520.6428
$ws.Range("I68").Value = 1960
$ws.Range("K68").Value = 1960
$ws.Range("M68").Value = -1211
$ws.Range("H71").Value = 2520.6428
$ws.Range("I71").Value = 1960
$ws.Range("K71").Value = 9800
$ws.Range("M71").Value = -6056
$ws.Range("H82").Value = 3501
$ws.Range("I82").Value = 1750
$ws.Range("K82").Value = 1750
$ws.Range("M82").Value = -1389
$ws.Range("H85").Value = 3501
$ws.Range("I85").Value = 1750
$ws.Range("K85").Value = 1750
$ws.Range("M85").Value = -502
$ws.Range("H122").Value = 1156563
$ws.Range("I122").Value = 1228442
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 3685326
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -3682876
$ws.Range("N122").Value = -24400
$ws.Range("H136").Value = 2473.5715
$ws.Range("I136").Value = 1908.9333
$ws.Range("J136").Value = 3885.1667
$ws.Range("K136").Value = 5726.7999
$ws.Range("L136").Value = 11655.5001
$ws.Range("M136").Value = -3176.7999
$ws.Range("N136").Value = -16755.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4999.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4375.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 24997.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21877.5
$ws.Range("N65").ClearContents()
